# Highlight the "listar clientes" and "listar productos" requirement
# paragraphs in yellow (wdYellow = 7), matching the target diff which
# adds <w:highlight w:val="yellow"/> to the runs of those two paragraphs.

$d = $word.ActiveDocument

$targets = @(
    "El sistema debe listar todos los clientes (panel principal);",
    "El sistema debe listar todos los productos (panel principal);"
)

foreach ($t in $targets) {
    $rng = $d.Content
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.HighlightColorIndex = 7
    }
}
